# Insert a new Apio price record for Macroferia Regional de Talca as row 26,
# pushing the existing rows 26..165 down to 27..166 (dimension grows to A1:R166).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 26..165 down by one, leaving a blank (but style-preserving) row 26.
$ws.Rows.Item(26).Insert()

# Populate the new row 26 with the inserted record.
$ws.Cells.Item(26, 1).Value = 5
$ws.Cells.Item(26, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(26, 3).Value = "Maule"
$ws.Cells.Item(26, 4).Value = 44687
$ws.Cells.Item(26, 5).Value = 7
$ws.Cells.Item(26, 6).Value = 100112017
$ws.Cells.Item(26, 7).Value = "Apio"
$ws.Cells.Item(26, 8).Value = "Americana (o)"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 500
$ws.Cells.Item(26, 11).Value = 7500
$ws.Cells.Item(26, 12).Value = 7500
$ws.Cells.Item(26, 13).Value = 7500
$ws.Cells.Item(26, 14).Value = "`$/docena de matas"
$ws.Cells.Item(26, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(26, 16).Value = 1250
$ws.Cells.Item(26, 17).Value = 6
$ws.Cells.Item(26, 18).Value = "Hortaliza"
